$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1928  # was 1923
$ws.Range("F6").Value = 1289  # was 1286
$ws.Range("F7").Value = 1289  # was 1286
$ws.Range("F9").Value = 1598  # was 1597
$ws.Range("F13").Value = 1679  # was 1677
$ws.Range("F15").Value = 1843  # was 1842
$ws.Range("F19").Value = 504  # was 502
$ws.Range("F20").Value = 1588  # was 1587
$ws.Range("F24").Value = 1091  # was 1090
$ws.Range("F25").Value = 2369  # was 2368
$ws.Range("F26").Value = 425  # was 424
$ws.Range("F28").Value = 1008  # was 715
$ws.Range("F29").Value = 4513  # was 4508
$ws.Range("F30").Value = 103  # was 101
$ws.Range("F31").Value = 31  # was 30
$ws.Range("F33").Value = 167  # was 165
$ws.Range("F35").Value = 1239  # was 1238
$ws.Range("F36").Value = 983  # was 982

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 68  # was 67
$ws.Range("F9").Value = 170  # was 169
$ws.Range("F20").Value = 182  # was 179
$ws.Range("F21").Value = 9  # was 7
$ws.Range("F22").Value = 9  # was 7
$ws.Range("F23").Value = 201  # was 200
$ws.Range("F34").Value = 468  # was 467
$ws.Range("F40").Value = 41  # was 38
$ws.Range("F47").Value = 36  # was 35

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2563  # was 2562
$ws.Range("F9").Value = 3081  # was 3076
$ws.Range("F10").Value = 600  # was 599
$ws.Range("F13").Value = 21  # was 10
$ws.Range("F14").Value = 21  # was 11
$ws.Range("F15").Value = 4  # was 2

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2563  # was 2562
$ws.Range("F8").Value = 3081  # was 3076
$ws.Range("F9").Value = 600  # was 599
$ws.Range("F12").Value = 1928  # was 1923
$ws.Range("F14").Value = 21  # was 10
$ws.Range("F15").Value = 21  # was 11
$ws.Range("F16").Value = 1289  # was 1286
$ws.Range("F17").Value = 68  # was 67
$ws.Range("F20").Value = 170  # was 169
$ws.Range("F21").Value = 1679  # was 1677
$ws.Range("F23").Value = 1843  # was 1842
$ws.Range("F26").Value = 504  # was 502
$ws.Range("F28").Value = 1588  # was 1587
$ws.Range("F30").Value = 182  # was 179
$ws.Range("F31").Value = 9  # was 7
$ws.Range("F33").Value = 201  # was 200
$ws.Range("F34").Value = 1091  # was 1090
$ws.Range("F37").Value = 425  # was 424
$ws.Range("F41").Value = 4513  # was 4508
$ws.Range("F42").Value = 31  # was 30
$ws.Range("F44").Value = 41  # was 38
$ws.Range("F46").Value = 167  # was 165
$ws.Range("F50").Value = 36  # was 35
$ws.Range("F51").Value = 1239  # was 1238
$ws.Range("F52").Value = 983  # was 982
